$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 37.6875
$ws.Range("I38").Value = 37.6875
$ws.Range("K38").Value = 113.0625
$ws.Range("M38").Value = 258.9375

$ws.Range("H51").Value = 3646.8125
$ws.Range("I51").Value = 1174.5
$ws.Range("J51").Value = 4000
$ws.Range("K51").Value = 1174.5
$ws.Range("L51").Value = 4000
$ws.Range("M51").Value = -690.5
$ws.Range("N51").Value = -4968

$ws.Range("H138").Value = 2375.8
$ws.Range("I138").Value = 1607.1666
$ws.Range("J138").Value = 2888.2222
$ws.Range("K138").Value = 4821.4998
$ws.Range("L138").Value = 8664.6666
$ws.Range("M138").Value = 318.5002000000004
$ws.Range("N138").Value = -18944.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 66000.31
$ws.Range("I61").Value = 3667.1667
$ws.Range("K61").Value = 3667.1667
$ws.Range("M61").Value = -3455.1667

$ws.Range("H132").Value = 2454.9656
$ws.Range("I132").Value = 2156.7727
$ws.Range("J132").Value = 3392.1428
$ws.Range("K132").Value = 6470.3181
$ws.Range("L132").Value = 10176.4284
$ws.Range("M132").Value = -3940.3181
$ws.Range("N132").Value = -15236.4284

$ws.Range("H136").Value = 66000.31
$ws.Range("I136").Value = 3667.1667
$ws.Range("K136").Value = 11001.5001
$ws.Range("M136").Value = -8451.500100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4703.1816
$ws.Range("I86").Value = 4593.25
$ws.Range("J86").Value = 4996.3335
$ws.Range("K86").Value = 4593.25
$ws.Range("L86").Value = 4996.3335
$ws.Range("M86").Value = -3470.25
$ws.Range("N86").Value = -7242.3335

$ws.Range("H89").Value = 4703.1816
$ws.Range("I89").Value = 4593.25
$ws.Range("J89").Value = 4996.3335
$ws.Range("K89").Value = 22966.25
$ws.Range("L89").Value = 24981.6675
$ws.Range("M89").Value = -17350.25
$ws.Range("N89").Value = -36213.6675

$ws.Range("H105").Value = 73980.28999999999
$ws.Range("I105").Value = 167875
$ws.Range("K105").Value = 167875
$ws.Range("M105").Value = -166128

$ws.Range("H134").Value = 4510.6333
$ws.Range("I134").Value = 2775.35
$ws.Range("J134").Value = 7981.2
$ws.Range("K134").Value = 8326.049999999999
$ws.Range("L134").Value = 23943.6
$ws.Range("M134").Value = -5791.049999999999
$ws.Range("N134").Value = -29013.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 8204.76
$ws.Range("I7").Value = 14349.714
$ws.Range("K7").Value = 14349.714
$ws.Range("M7").Value = -14236.714

$ws.Range("H31").Value = 3533.524
$ws.Range("I31").Value = 2194.2856
$ws.Range("J31").Value = 6212
$ws.Range("K31").Value = 2194.2856
$ws.Range("L31").Value = 6212
$ws.Range("M31").Value = -1899.2856
$ws.Range("N31").Value = -6802

$ws.Range("H34").Value = 3533.524
$ws.Range("I34").Value = 2194.2856
$ws.Range("J34").Value = 6212
$ws.Range("K34").Value = 2194.2856
$ws.Range("L34").Value = 6212
$ws.Range("M34").Value = -1992.2856
$ws.Range("N34").Value = -6616

$ws.Range("H41").Value = 13363.637

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H94").Value = 949
$ws.Range("I94").Value = 1000
$ws.Range("J94").Value = 898
$ws.Range("K94").Value = 1000
$ws.Range("L94").Value = 898
$ws.Range("M94").Value = -549
$ws.Range("N94").Value = -1800

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6186.625
$ws.Range("I56").Value = 6186.625
$ws.Range("K56").Value = 6186.625
$ws.Range("M56").Value = -5656.625

$ws.Range("H70").Value = 600
$ws.Range("J70").Value = 600
$ws.Range("L70").Value = 1800
$ws.Range("N70").Value = -2430

$ws.Range("H73").Value = 600
$ws.Range("J73").Value = 600
$ws.Range("L73").Value = 1800
$ws.Range("N73").Value = -3984

$ws.Range("H121").Value = 1152.5333
$ws.Range("J121").Value = 1808.5714
$ws.Range("L121").Value = 5425.7142
$ws.Range("N121").Value = -8045.7142

$ws.Range("H138").Value = 6724.778
$ws.Range("I138").Value = 8524
$ws.Range("K138").Value = 25572
$ws.Range("M138").Value = -20432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 10007043
$ws.Range("J33").Value = 10007483
$ws.Range("L33").Value = 10007483
$ws.Range("N33").Value = -10007987

$ws.Range("H52").Value = 16413.285
$ws.Range("J52").Value = 15815.5
$ws.Range("L52").Value = 15815.5
$ws.Range("N52").Value = -16333.5

$ws.Range("H80").Value = 905.5
$ws.Range("I80").Value = 917.6667
$ws.Range("J80").Value = 869
$ws.Range("K80").Value = 917.6667
$ws.Range("L80").Value = 869
$ws.Range("M80").Value = 80.33330000000001
$ws.Range("N80").Value = -2865

$ws.Range("H83").Value = 905.5
$ws.Range("I83").Value = 917.6667
$ws.Range("J83").Value = 869
$ws.Range("K83").Value = 4588.3335
$ws.Range("L83").Value = 4345
$ws.Range("M83").Value = 403.6665000000003
$ws.Range("N83").Value = -14329

$ws.Range("H122").Value = 12226798
$ws.Range("I122").Value = 13754835
$ws.Range("K122").Value = 41264505
$ws.Range("M122").Value = -41262055

$ws.Range("H132").Value = 4190.3335
$ws.Range("I132").Value = 2793.7
$ws.Range("K132").Value = 8381.099999999999
$ws.Range("M132").Value = -5851.099999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3190.8333
$ws.Range("I7").Value = 2242.8572
$ws.Range("K7").Value = 2242.8572
$ws.Range("M7").Value = -2130.8572

$ws.Range("H46").Value = 7793.1875
$ws.Range("I46").Value = 8649.714
$ws.Range("J46").Value = 1797.5
$ws.Range("K46").Value = 8649.714
$ws.Range("L46").Value = 1797.5
$ws.Range("M46").Value = -8461.714
$ws.Range("N46").Value = -2173.5

$ws.Range("H75").Value = 22222
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 22222
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H122").Value = 118185576
$ws.Range("I122").Value = 250003630
$ws.Range("J122").Value = 42860976
$ws.Range("K122").Value = 750010890
$ws.Range("L122").Value = 128582928
$ws.Range("M122").Value = -750008440
$ws.Range("N122").Value = -128587828

$ws.Range("H126").Value = 3190.8333
$ws.Range("I126").Value = 2242.8572
$ws.Range("K126").Value = 6728.571599999999
$ws.Range("M126").Value = -4258.571599999999

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws.Range("H132").Value = 8400
$ws.Range("I132").Value = 8466.666999999999
$ws.Range("K132").Value = 25400.001
$ws.Range("M132").Value = -22870.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 21887
$ws.Range("I22").Value = 29999
$ws.Range("J22").Value = 13775
$ws.Range("K22").Value = 29999
$ws.Range("L22").Value = 13775
$ws.Range("M22").Value = -29706
$ws.Range("N22").Value = -14361

$ws.Range("H52").Value = 22598
$ws.Range("I52").Value = 24497.5
$ws.Range("J52").Value = 15000
$ws.Range("K52").Value = 24497.5
$ws.Range("L52").Value = 15000
$ws.Range("M52").Value = -24271.5
$ws.Range("N52").Value = -15452

$ws.Range("H126").Value = 5199.7334
$ws.Range("I126").Value = 4925.8237
$ws.Range("J126").Value = 5557.923
$ws.Range("K126").Value = 14777.4711
$ws.Range("L126").Value = 16673.769
$ws.Range("M126").Value = -12307.4711
$ws.Range("N126").Value = -21613.769

$ws.Range("H132").Value = 2060.3333
$ws.Range("I132").Value = 1823.4546
$ws.Range("J132").Value = 2432.5715
$ws.Range("K132").Value = 5470.3638
$ws.Range("L132").Value = 7297.7145
$ws.Range("M132").Value = -2940.3638
$ws.Range("N132").Value = -12357.7145

$ws.Range("H136").Value = 899.8
$ws.Range("I136").Value = 999.75
$ws.Range("K136").Value = 2999.25
$ws.Range("M136").Value = -449.25
